# distance_comparisons.xlsx - "Changed the stretching to y-axis (rotated)"
#
# The anisotropy stretching factor ($C$18 / $C$19) that used to be applied
# to the x-distance column (L) is moved to the y-distance column (M)
# instead, for both the first group of points (rows 3-8, using $C$18/$D$18)
# and the second group (rows 9-14, using $C$19/$D$19).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- First block (rows 3-8), anchored on row 3 / $C$18 ---
$ws.Range("L3").Formula = '=((I$3-I3)^2)'
$ws.Range("M3").Formula = '=$C$18*((J$3-J3)^2)'

$ws.Range("L4:L8").Formula = '=((I$3-I4)^2)'
$ws.Range("M4:M8").Formula = '=$C$18*((J$3-J4)^2)'

# --- Second block (rows 9-14), anchored on row 9 / $C$19 ---
$ws.Range("L9").Formula = '=((I$9-I9)^2)'
$ws.Range("M9").Formula = '=$C$19*((J$9-J9)^2)'

$ws.Range("L10:L14").Formula = '=((I$9-I10)^2)'
$ws.Range("M10:M14").Formula = '=$C$19*((J$9-J10)^2)'

# --- Distance column recomputed as one contiguous block (unchanged text,
#     just re-entered so the whole range now shares one formula group) ---
$ws.Range("O3:O14").Formula = '=SQRT(SUM(L3:N3))'

# --- "Spreadsheet results" tables: hard-coded copies of the O column,
#     re-pasted with the new values and reformatted from the custom
#     7-decimal scientific format to the built-in 0.00E+00 format. ---
$ws.Range("C23").Value = 1.266
$ws.Range("D23").Value = 2.532
$ws.Range("E23").Value = 4.939
$ws.Range("F23").Value = 4.637
$ws.Range("G23").Value = 4.67

$ws.Range("C25").Value = 1
$ws.Range("D25").Value = 2
$ws.Range("E25").Value = 3.873
$ws.Range("F25").Value = 4
$ws.Range("G25").Value = 4.359

$ws.Range("C23:G23").NumberFormat = "0.00E+00"
$ws.Range("C25:G25").NumberFormat = "0.00E+00"

# --- Selection moved to P35 ---
$ws.Range("P35").Select()
